$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.821.83"
$ws.Range("E2").Value = "  -2.21%  "

$ws.Range("D3").Value = "3.940.15"

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").Value = "'533.39"
$ws.Range("E5").Value = "  +2.24%  "

$ws.Range("D6").Value = "'147.70"
$ws.Range("E6").Value = "  -0.14%  "

$ws.Range("D7").Value = "3.934.47"
$ws.Range("E7").Value = "  -2.74%  "

$ws.Range("D8").Value = "'0.686"
$ws.Range("E8").Value = "  -3.94%  "

$ws.Range("D9").Value = "'0.999"
$ws.Range("E9").Value = "  -0.06%  "

$ws.Range("D10").Value = "'0.736"
$ws.Range("E10").Value = "  -4.70%  "

$ws.Range("D11").Value = "'0.165"
$ws.Range("E11").Value = "  -8.62%  "

$ws.Range("D12").Value = "'55.16"
$ws.Range("E12").Value = "  +13.88%  "

$ws.Range("D13").Value = "'0.0000315"
$ws.Range("E13").Value = "  -5.99%  "

$ws.Range("D14").Value = "'10.56"
$ws.Range("E14").Value = "  -5.56%  "

$ws.Range("D15").Value = "4.560.72"
$ws.Range("E15").Value = "  -2.72%  "

$ws.Range("D16").Value = "3.937.81"
$ws.Range("E16").Value = "  -3.04%  "

$ws.Range("D17").Value = "'20.48"
$ws.Range("E17").Value = "  -3.98%  "

$ws.Range("D18").Value = "'13.81"
$ws.Range("E18").Value = "  -3.46%  "

$ws.Range("E19").Value = "  -1.68%  "

$ws.Range("E20").Value = "  -5.04%  "

$ws.Range("D21").Value = "70.687.47"
$ws.Range("E21").Value = "  -2.28%  "

$ws.Range("D22").Value = "'421.64"
$ws.Range("E22").Value = "  -5.48%  "

$ws.Range("D23").Value = "'3.59"
$ws.Range("E23").Value = "  -0.37%  "

$ws.Range("D24").Value = "'97.05"
$ws.Range("E24").Value = "  -7.58%  "

$ws.Range("D25").Value = "'4.20"
$ws.Range("E25").Value = "  +3.84%  "

$ws.Range("D26").Value = "'14.40"
$ws.Range("E26").Value = "  -5.69%  "

$ws.Range("D27").Value = "'11.31"

$ws.Range("D28").Value = "'3.80"
$ws.Range("E28").Value = "  +15.10%  "

$ws.Range("D29").Value = "'10.61"
$ws.Range("E29").Value = "  -4.70%  "

$ws.Range("D30").Value = "'5.86"
$ws.Range("E30").Value = "  +0.83%  "

$ws.Range("D31").Value = "'36.29"
$ws.Range("E31").Value = "  -4.56%  "

$ws.Range("D32").Value = "'7.82"
$ws.Range("E32").Value = "  +15.71%  "

$ws.Range("D33").Value = "'50.79"
$ws.Range("E33").Value = "  +19.35%  "

$ws.Range("D34").Value = "'0.131"
$ws.Range("E34").Value = "  -0.18%  "

$ws.Range("D35").Value = "'13.31"
$ws.Range("E35").Value = "  -3.58%  "

$ws.Range("D36").Value = "'683.53"
$ws.Range("E36").Value = "  +0.18%  "

$ws.Range("D37").Value = "'65.28"
$ws.Range("E37").Value = "  -3.97%  "

$ws.Range("D38").Value = "'0.437"
$ws.Range("E38").Value = "  +1.33%  "

$ws.Range("B39").Value = "PEPE"
$ws.Range("C39").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D39").Value = "0.0₃0811"
$ws.Range("E39").Value = "  -7.36%  "

$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "'0.148"
$ws.Range("E40").Value = "  -2.91%  "

$ws.Range("D41").Value = "'3.36"
$ws.Range("E41").Value = "  -3.98%  "

$ws.Range("E42").Value = "  -0.05%  "

$ws.Range("D43").Value = "'0.999"
$ws.Range("E43").Value = "  +0.15%  "

$ws.Range("D44").Value = "'0.0479"
$ws.Range("E44").Value = "  -4.82%  "

$ws.Range("E45").Value = "  -0.76%  "

$ws.Range("E46").Value = "  -5.53%  "

$ws.Range("E47").Value = "  +2.95%  "

$ws.Range("E48").Value = "  -1.83%  "

$ws.Range("D49").Value = "'3.32"
$ws.Range("E49").Value = "  -5.17%  "

$ws.Range("D50").Value = "'2.99"
$ws.Range("E50").Value = "  -3.20%  "

$ws.Range("D51").Value = "'144.77"
$ws.Range("E51").Value = "  -0.20%  "
